$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 4,23
$arr[0,0] = 0.00123571207908557
$arr[0,1] = 0.00046339202965709
$arr[0,2] = 0.949799196787149
$arr[0,3] = 0.00108124806919988
$arr[0,4] = 0.999073215940686
$arr[0,5] = 0.965091133765833
$arr[0,6] = 0.00494284831634229
$arr[0,7] = 0.000617856039542787
$arr[0,8] = 0.98671609514983
$arr[0,9] = 0.00417052826691381
$arr[0,10] = 0.997683039851715
$arr[0,11] = 0.000154464009885697
$arr[0,12] = 0.00123571207908557
$arr[0,13] = 0.00370713623725672
$arr[0,14] = 0
$arr[0,15] = 0.994748223663886
$arr[0,16] = 0.00278035217794254
$arr[0,17] = 0
$arr[0,18] = 0.000154464009885697
$arr[0,19] = 0.00046339202965709
$arr[0,20] = 0.993049119555144
$arr[0,21] = 0.00046339202965709
$arr[0,22] = 0.0358356502934816
$arr[1,0] = 0
$arr[1,1] = 0
$arr[1,2] = 0.00169910410874266
$arr[1,3] = 0.998300895891257
$arr[1,4] = 0.00046339202965709
$arr[1,5] = 0.000308928019771393
$arr[1,6] = 0
$arr[1,7] = 0.00139017608897127
$arr[1,8] = 0.00046339202965709
$arr[1,9] = 0
$arr[1,10] = 0.00046339202965709
$arr[1,11] = 0.938214396045721
$arr[1,12] = 0.00123571207908557
$arr[1,13] = 0
$arr[1,14] = 0.0607043558850788
$arr[1,15] = 0.000154464009885697
$arr[1,16] = 0.000154464009885697
$arr[1,17] = 0.0531356194006796
$arr[1,18] = 0.999691071980229
$arr[1,19] = 0.999382143960457
$arr[1,20] = 0.000772320049428483
$arr[1,21] = 0.000308928019771393
$arr[1,22] = 0.000772320049428483
$arr[2,0] = 0.998455359901143
$arr[2,1] = 0.999227679950571
$arr[2,2] = 0.0421686746987952
$arr[2,3] = 0
$arr[2,4] = 0.000154464009885697
$arr[2,5] = 0.0287303058387396
$arr[2,6] = 0.993821439604572
$arr[2,7] = 0
$arr[2,8] = 0.00880444856348471
$arr[2,9] = 0.995675007723201
$arr[2,10] = 0.00169910410874266
$arr[2,11] = 0
$arr[2,12] = 0.000308928019771393
$arr[2,13] = 0.995983935742972
$arr[2,14] = 0.000154464009885697
$arr[2,15] = 0.00417052826691381
$arr[2,16] = 0.996138399752858
$arr[2,17] = 0.000154464009885697
$arr[2,18] = 0
$arr[2,19] = 0
$arr[2,20] = 0.00200803212851406
$arr[2,21] = 0.999073215940686
$arr[2,22] = 0.962156317578004
$arr[3,0] = 0
$arr[3,1] = 0.000308928019771393
$arr[3,2] = 0.00169910410874266
$arr[3,3] = 0.000308928019771393
$arr[3,4] = 0
$arr[3,5] = 0.000772320049428483
$arr[3,6] = 0.000154464009885697
$arr[3,7] = 0.997991967871486
$arr[3,8] = 0.000772320049428483
$arr[3,9] = 0
$arr[3,10] = 0
$arr[3,11] = 0.0529811553907939
$arr[3,12] = 0.996601791782515
$arr[3,13] = 0.000154464009885697
$arr[3,14] = 0.934198331788693
$arr[3,15] = 0
$arr[3,16] = 0.000617856039542787
$arr[3,17] = 0.939759036144578
$arr[3,18] = 0.000154464009885697
$arr[3,19] = 0
$arr[3,20] = 0.00278035217794254
$arr[3,21] = 0.000154464009885697
$arr[3,22] = 0

$ws.Range("B2:X5").Value = $arr
